$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows for column I (I0) and column J (IF)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 8

$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 8
